$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "Description" column (B) for the metadata rows that previously
# only had a "Variable name" (column A) and no description, plus refine the
# wording for a few related rows. Values are written in the same order the
# author filled them in (rows 54-63, then 42-53) so that the shared-strings
# table gets appended in the same sequence as the target workbook.
$ws.Range('B54').Value = 'encoded unique identifier, joining stock ID and management program name'
$ws.Range('B55').Value = 'encoded RAM region containing the majority of the spatial distribution of the stock'
$ws.Range('B56').Value = 'encoded International Standard Statistical Classification of Aquatic Animals and Plants classification'
$ws.Range('B57').Value = 'encoded FAO region containing the majority of the spatial distribution of the stock'
$ws.Range('B58').Value = 'number of containing EEZs for this stock - 1, at a 0% threshold'
$ws.Range('B59').Value = 'number of containing EEZs for this stock - 1, at a 5% threshold'
$ws.Range('B60').Value = 'number of containing EEZs for this stock - 1, at a 10% threshold'
$ws.Range('B61').Value = 'number of containing EEZs for this stock - 1, at a 15% threshold'
$ws.Range('B62').Value = 'number of containing EEZs for this stock - 1, at a 20% threshold'
$ws.Range('B63').Value = 'number of containing EEZs for this stock - 1, at a 25% threshold'
$ws.Range('B42').Value = 'binary indicator of whether or not a stock crosses two or more jurisdictions, at a 0% threshold'
$ws.Range('B43').Value = 'counter of total stocks shared, at a 0% threshold'
$ws.Range('B44').Value = 'binary indicator of whether or not a stock crosses two or more jurisdictions, at a 5% threshold'
$ws.Range('B45').Value = 'counter of total stocks shared, at a 5% threshold'
$ws.Range('B46').Value = 'binary indicator of whether or not a stock crosses two or more jurisdictions, at a 10% threshold'
$ws.Range('B47').Value = 'counter of total stocks shared, at a 10% threshold'
$ws.Range('B48').Value = 'binary indicator of whether or not a stock crosses two or more jurisdictions, at a 15% threshold'
$ws.Range('B49').Value = 'counter of total stocks shared, at a 15% threshold'
$ws.Range('B50').Value = 'binary indicator of whether or not a stock crosses two or more jurisdictions, at a 20% threshold'
$ws.Range('B51').Value = 'counter of total stocks shared, at a 20% threshold'
$ws.Range('B52').Value = 'binary indicator of whether or not a stock crosses two or more jurisdictions, at a 25% threshold'
$ws.Range('B53').Value = 'counter of total stocks shared, at a 25% threshold'

# Re-select cell B54, matching the final cursor position left by the author.
$ws.Range('B54').Select()

# Best-effort match of the column widths recorded by the author's Excel
# session (this engine only persists ColumnWidth, quantized to 1/6 of a
# character, so we pick the input that lands nearest the recorded widths).
$ws.Columns.Item(1).ColumnWidth = 24
$ws.Columns.Item(2).ColumnWidth = 10.3
